$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "37.288.98"
Set-TextValue "E2" "  +5.27%  "
Set-TextValue "D3" "1.926.96"
Set-TextValue "E3" "  +2.00%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "249.62"
Set-TextValue "E5" "  +1.22%  "
Set-TextValue "D6" "0.688"
Set-TextValue "E6" "  -0.65%  "
Set-TextValue "E7" "  -0.07%  "
Set-TextValue "D8" "47.48"
Set-TextValue "E8" "  +9.90%  "
Set-TextValue "E9" "  +5.81%  "
Set-TextValue "D10" "58.21"
Set-TextValue "E10" "  +6.25%  "
Set-TextValue "D11" "0.0764"
Set-TextValue "E11" "  +2.42%  "
Set-TextValue "E12" "  +1.51%  "
Set-TextValue "D13" "15.62"
Set-TextValue "D14" "0.825"
Set-TextValue "E14" "  +6.77%  "
Set-TextValue "D15" "2.202.78"
Set-TextValue "E15" "  +1.86%  "
Set-TextValue "D16" "5.14"
Set-TextValue "E16" "  +2.11%  "
Set-TextValue "D17" "1.922.03"
Set-TextValue "E17" "  +1.76%  "
Set-TextValue "D18" "37.219.21"
Set-TextValue "E18" "  +5.11%  "
Set-TextValue "D19" "74.84"
Set-TextValue "E19" "  +1.72%  "
Set-TextValue "D20" "0.0₃0859"
Set-TextValue "E20" "  +3.88%  "
Set-TextValue "D21" "13.68"
Set-TextValue "E21" "  +6.43%  "
Set-TextValue "D22" "251.17"
Set-TextValue "E22" "  +2.32%  "
Set-TextValue "D23" "5.19"
Set-TextValue "E23" "  +0.68%  "
Set-TextValue "E24" "  -0.08%  "
Set-TextValue "D25" "2.52"
Set-TextValue "E25" "  -4.74%  "
Set-TextValue "D26" "168.47"
Set-TextValue "E26" "  +1.23%  "
Set-TextValue "D27" "2.11"
Set-TextValue "E27" "  -2.43%  "
Set-TextValue "D28" "8.81"
Set-TextValue "E28" "  +2.13%  "
Set-TextValue "D29" "18.78"
Set-TextValue "E29" "  +2.49%  "
Set-TextValue "E30" "  +0.04%  "
Set-TextValue "D31" "4.57"
Set-TextValue "E31" "  +6.09%  "
Set-TextValue "D32" "0.0611"
Set-TextValue "E32" "  +2.36%  "
Set-TextValue "D33" "0.0917"
Set-TextValue "E34" "  +2.19%  "
Set-TextValue "D35" "1.89"
Set-TextValue "E35" "  +0.03%  "
Set-TextValue "E36" "  +0.09%  "
Set-TextValue "D37" "19.14"
Set-TextValue "E37" "  +39.15%  "
Set-TextValue "D38" "0.898"
Set-TextValue "E38" "  +4.63%  "
Set-TextValue "D39" "1.44"
Set-TextValue "E39" "  -1.35%  "
Set-TextValue "E40" "  -0.26%  "
Set-TextValue "D41" "104.96"
Set-TextValue "E41" "  +6.96%  "
Set-TextValue "E42" "  +2.36%  "
Set-TextValue "E43" "  +1.50%  "
Set-TextValue "D44" "2.94"
Set-TextValue "E44" "  +22.57%  "
Set-TextValue "D45" "1.10"
Set-TextValue "E45" "  +1.72%  "
Set-TextValue "D46" "1.346.54"
Set-TextValue "E46" "  +1.20%  "
Set-TextValue "D47" "2.40"
Set-TextValue "E47" "  +0.06%  "
Set-TextValue "D48" "0.0838"
Set-TextValue "E48" "  +3.51%  "
Set-TextValue "D49" "2.81"
Set-TextValue "E49" "  +2.46%  "
Set-TextValue "D50" "6.39"
Set-TextValue "E50" "  +1.30%  "
Set-TextValue "D51" "3.77"
Set-TextValue "E51" "  +13.08%  "
